$wb = $excel.ActiveWorkbook

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(111, 8).Value = 1458.1818
$ws.Cells.Item(111, 9).Value = 1276.3334
$ws.Cells.Item(111, 10).Value = 1676.4
$ws.Cells.Item(111, 11).Value = 3829.0002
$ws.Cells.Item(111, 12).Value = 5029.200000000001
$ws.Cells.Item(111, 13).Value = -762.0001999999999
$ws.Cells.Item(111, 14).Value = -11163.2

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3381.9062
$ws.Cells.Item(132, 9).Value = 3470.077
$ws.Cells.Item(132, 10).Value = 2999.8333
$ws.Cells.Item(132, 11).Value = 10410.231
$ws.Cells.Item(132, 12).Value = 8999.499899999999
$ws.Cells.Item(132, 13).Value = -7880.231
$ws.Cells.Item(132, 14).Value = -14059.4999

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 183.875
$ws.Cells.Item(4, 9).Value = 189.6
$ws.Cells.Item(4, 10).Value = 174.33333
$ws.Cells.Item(4, 11).Value = 189.6
$ws.Cells.Item(4, 12).Value = 174.33333
$ws.Cells.Item(4, 13).Value = -73.59999999999999
$ws.Cells.Item(4, 14).Value = -406.33333

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 145371.42
$ws.Cells.Item(122, 9).Value = 251950
$ws.Cells.Item(122, 10).Value = 3266.6667
$ws.Cells.Item(122, 11).Value = 755850
$ws.Cells.Item(122, 12).Value = 9800.000100000001
$ws.Cells.Item(122, 13).Value = -753400
$ws.Cells.Item(122, 14).Value = -14700.0001

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1641218
$ws.Cells.Item(132, 9).Value = 4701.6523
$ws.Cells.Item(132, 10).Value = 3209546.2
$ws.Cells.Item(132, 11).Value = 14104.9569
$ws.Cells.Item(132, 12).Value = 9628638.600000001
$ws.Cells.Item(132, 13).Value = -11574.9569
$ws.Cells.Item(132, 14).Value = -9633698.600000001

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 300
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 12).Value = 300
$ws.Cells.Item(22, 14).Value = -646

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 93023.91
$ws.Cells.Item(107, 9).Value = 335100
$ws.Cells.Item(107, 11).Value = 335100
$ws.Cells.Item(107, 13).Value = -333180

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2793.257
$ws.Cells.Item(134, 9).Value = 2877.9167
$ws.Cells.Item(134, 10).Value = 2608.5454
$ws.Cells.Item(134, 11).Value = 8633.750100000001
$ws.Cells.Item(134, 12).Value = 7825.6362
$ws.Cells.Item(134, 13).Value = -6098.750100000001
$ws.Cells.Item(134, 14).Value = -12895.6362

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5534.75
$ws.Cells.Item(31, 9).Value = 1736.1
$ws.Cells.Item(31, 11).Value = 1736.1
$ws.Cells.Item(31, 13).Value = -1441.1

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5534.75
$ws.Cells.Item(34, 9).Value = 1736.1
$ws.Cells.Item(34, 11).Value = 1736.1
$ws.Cells.Item(34, 13).Value = -1534.1

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 3133.3333
$ws.Cells.Item(86, 9).Value = 3200
$ws.Cells.Item(86, 10).Value = 2800
$ws.Cells.Item(86, 11).Value = 3200
$ws.Cells.Item(86, 12).Value = 2800
$ws.Cells.Item(86, 13).Value = -2077
$ws.Cells.Item(86, 14).Value = -5046

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 3133.3333
$ws.Cells.Item(89, 9).Value = 3200
$ws.Cells.Item(89, 10).Value = 2800
$ws.Cells.Item(89, 11).Value = 16000
$ws.Cells.Item(89, 12).Value = 14000
$ws.Cells.Item(89, 13).Value = -10384
$ws.Cells.Item(89, 14).Value = -25232

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 846.6539
$ws.Cells.Item(107, 9).Value = 683.75
$ws.Cells.Item(107, 10).Value = 919.05554
$ws.Cells.Item(107, 11).Value = 683.75
$ws.Cells.Item(107, 12).Value = 919.05554
$ws.Cells.Item(107, 13).Value = 1236.25
$ws.Cells.Item(107, 14).Value = -4759.05554

# CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 1000
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()

# CUL row 31
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(31, 8).Value = 2199.3125
$ws.Cells.Item(31, 9).Value = 866.3333
$ws.Cells.Item(31, 10).Value = 2506.923
$ws.Cells.Item(31, 11).Value = 2598.9999
$ws.Cells.Item(31, 12).Value = 7520.768999999999
$ws.Cells.Item(31, 13).Value = -2310.9999
$ws.Cells.Item(31, 14).Value = -8096.768999999999

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 11114245
$ws.Cells.Item(32, 10).Value = 11114245
$ws.Cells.Item(32, 12).Value = 33342735
$ws.Cells.Item(32, 14).Value = -33343301

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 87789
$ws.Cells.Item(137, 9).Value = 50914.5
$ws.Cells.Item(137, 11).Value = 152743.5
$ws.Cells.Item(137, 13).Value = -147643.5

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 1068.6364
$ws.Cells.Item(140, 9).Value = 935.5
$ws.Cells.Item(140, 10).Value = 2400
$ws.Cells.Item(140, 11).Value = 2806.5
$ws.Cells.Item(140, 12).Value = 7200
$ws.Cells.Item(140, 13).Value = 2373.5
$ws.Cells.Item(140, 14).Value = -17560

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2283.7144
$ws.Cells.Item(102, 9).Value = 2235.2307
$ws.Cells.Item(102, 10).Value = 2362.5
$ws.Cells.Item(102, 11).Value = 2235.2307
$ws.Cells.Item(102, 12).Value = 2362.5
$ws.Cells.Item(102, 13).Value = -613.2307000000001
$ws.Cells.Item(102, 14).Value = -5606.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 103405.5
$ws.Cells.Item(113, 9).Value = 146436.42
$ws.Cells.Item(113, 11).Value = 146436.42
$ws.Cells.Item(113, 13).Value = -144266.42

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1956.5454
$ws.Cells.Item(122, 9).Value = 2056.2222
$ws.Cells.Item(122, 10).Value = 1508
$ws.Cells.Item(122, 11).Value = 6168.6666
$ws.Cells.Item(122, 12).Value = 4524
$ws.Cells.Item(122, 13).Value = -3718.6666
$ws.Cells.Item(122, 14).Value = -9424

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 15947.143
$ws.Cells.Item(123, 10).Value = 15947.143
$ws.Cells.Item(123, 12).Value = 15947.143
$ws.Cells.Item(123, 14).Value = -20847.143

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 997.25
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 997.25
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 2991.75
$ws.Cells.Item(126, 14).Value = -7931.75
$ws.Cells.Item(126, 13).ClearContents()

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4502
$ws.Cells.Item(132, 9).Value = 4057.6667
$ws.Cells.Item(132, 10).Value = 5168.5
$ws.Cells.Item(132, 11).Value = 12173.0001
$ws.Cells.Item(132, 12).Value = 15505.5
$ws.Cells.Item(132, 13).Value = -9643.000100000001
$ws.Cells.Item(132, 14).Value = -20565.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4060.8572
$ws.Cells.Item(122, 9).Value = 3450.5881
$ws.Cells.Item(122, 10).Value = 4637.222
$ws.Cells.Item(122, 11).Value = 10351.7643
$ws.Cells.Item(122, 12).Value = 13911.666
$ws.Cells.Item(122, 13).Value = -7901.764299999999
$ws.Cells.Item(122, 14).Value = -18811.666

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2592.7856
$ws.Cells.Item(136, 9).Value = 2850.5
$ws.Cells.Item(136, 10).Value = 2399.5
$ws.Cells.Item(136, 11).Value = 8551.5
$ws.Cells.Item(136, 12).Value = 7198.5
$ws.Cells.Item(136, 13).Value = -6001.5
$ws.Cells.Item(136, 14).Value = -12298.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 8840887
$ws.Cells.Item(132, 9).Value = 3165
$ws.Cells.Item(132, 10).Value = 18230966
$ws.Cells.Item(132, 11).Value = 9495
$ws.Cells.Item(132, 12).Value = 54692898
$ws.Cells.Item(132, 13).Value = -6965
$ws.Cells.Item(132, 14).Value = -54697958
